{"js": "// Remove the empty paragraph that immediately follows the paragraph\n// ending in \" del paciente.\" (the \"Aporte\" section right before the\n// \"\u00bfQUE ES GESTI\u00d3N DE SERVICIO DE ATENCI\u00d3N?\" heading).\nconst body = context.document.body;\n\n// The phrase \" del paciente.\" (with leading space, from the run\n// xml:space=\"preserve\"> del paciente.</w:t>) is unique in the document,\n// so searching for it unambiguously locates the target paragraph.\nconst results = body.search(\" del paciente.\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the text \" del paciente.\" in the document.');\n}\n\n// Paragraph containing the matched text.\nconst targetPara = results.items[0].paragraphs.getFirst();\n\n// The paragraph immediately following it \u2014 this is the empty paragraph\n// that the diff removes.\nconst emptyPara = targetPara.getNext();\nemptyPara.load(\"text\");\nawait context.sync();\n\nif (emptyPara.text.trim() === \"\") {\n  emptyPara.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the empty paragraph that immediately follows the paragraph\n# ending in \" del paciente.\" (the \"Aporte\" section right before the\n# \"\u00bfQUE ES GESTI\u00d3N DE SERVICIO DE ATENCI\u00d3N?\" heading).\n$d = $word.ActiveDocument\n\n# The phrase \"del paciente.\" is unique in the document (only the\n# paragraph about \"...todo el expediente medico del paciente.\" contains\n# it), so matching on it unambiguously locates the target paragraph.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*del paciente.*\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -eq $null) {\n    throw \"Could not find the paragraph containing 'del paciente.'\"\n}\n\n# The paragraph immediately following it -- this is the empty paragraph\n# that the diff removes.\n$emptyPara = $targetPara.Next()\n\nif ($emptyPara.Range.Text.Trim().Length -eq 0) {\n    $emptyPara.Range.Delete()\n}\n"}
